$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.073.18'
$ws.Range('E2').Value = '  +4.61%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.239.72'
$ws.Range('E3').Value = '  +4.65%  '

# Row 4
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.26'
$ws.Range('E5').Value = '  +6.57%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('E6').Value = '  +2.53%  '

# Row 7
$ws.Range('E7').Value = '  +9.00%  '

# Row 8
$ws.Range('E8').Value = '  -0.18%  '

# Row 9
$ws.Range('E9').Value = '  +6.18%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.91'
$ws.Range('E10').Value = '  +7.14%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  +4.39%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.90'
$ws.Range('E12').Value = '  +4.90%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.102'
$ws.Range('E13').Value = '  +2.22%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.575.82'
$ws.Range('E14').Value = '  +4.55%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.64'
$ws.Range('E15').Value = '  +2.28%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.273.99'
$ws.Range('E16').Value = '  +6.55%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.791'
$ws.Range('E17').Value = '  +2.25%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.983.07'
$ws.Range('E18').Value = '  +4.72%  '

# Row 19
$ws.Range('E19').Value = '  +6.04%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.14'
$ws.Range('E20').Value = '  +2.91%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.98'
$ws.Range('E21').Value = '  +5.02%  '

# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.66'
$ws.Range('E22').Value = '  +2.77%  '

# Row 23
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.20'
$ws.Range('E23').Value = '  +17.04%  '

# Row 24
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.70'
$ws.Range('E24').Value = '  +2.73%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.77'
$ws.Range('E26').Value = '  +2.43%  '

# Row 27
$ws.Range('E27').Value = '  +2.14%  '

# Row 28
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.53'
$ws.Range('E28').Value = '  +28.23%  '

# Row 29
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  +5.29%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.24'
$ws.Range('E30').Value = '  +4.72%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.53'
$ws.Range('E31').Value = '  +1.89%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.21'
$ws.Range('E32').Value = '  +3.76%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0799'
$ws.Range('E33').Value = '  +6.73%  '

# Row 34
$ws.Range('E34').Value = '  +4.59%  '

# Row 35
$ws.Range('E35').Value = '  +2.26%  '

# Row 36
$ws.Range('E36').Value = '  +10.13%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.48'
$ws.Range('E37').Value = '  +10.43%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0329'
$ws.Range('E38').Value = '  +17.99%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.98'
$ws.Range('E39').Value = '  +12.40%  '

# Row 40
$ws.Range('E40').Value = '  +4.38%  '

# Row 41
$ws.Range('E41').Value = '  +11.04%  '

# Row 42
$ws.Range('E42').Value = '  +3.78%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '59.43'
$ws.Range('E43').Value = '  +4.49%  '

# Row 44
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.64'
$ws.Range('E44').Value = '  +6.33%  '

# Row 45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.76'
$ws.Range('E45').Value = '  +7.69%  '

# Row 46
$ws.Range('B46').Value = 'WOONetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.484'
$ws.Range('E46').Value = '  +30.78%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0990'
$ws.Range('E47').Value = '  +4.73%  '

# Row 48
$ws.Range('E48').Value = '  +13.71%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('E49').Value = '  +4.09%  '

# Row 50
$ws.Range('E50').Value = '  +4.78%  '

# Row 51
$ws.Range('E51').Value = '  +3.37%  '
